$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.945.17"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.816.71"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'310.13"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4649"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.3709"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "'0.07375"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'0.8736"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "1.788.59"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").Value = "'5.359"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "'6.519"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "'0.07058"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'91.63"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'0.000008737"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'14.74"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "26.968.00"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").Value = "'5.326"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "2.060.77"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "'1.910"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'151.84"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'2.155"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "'5.329"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'115.92"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "'0.08918"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'0.7587"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").Value = "'1.159"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").Value = "'4.482"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "'2.921"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'1.097"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'0.01960"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "'0.05265"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "'2.427"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").Value = "'2.936"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.261"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5352"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "'8.456"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "'0.4968"
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("D47").Value = "'10.36"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "'1.682"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'103.39"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'0.06295"
$ws.Range("E51").Value = "  -0.81%  "
